# Refresh the cryptocurrency price/volume snapshot (Price = column D,
# Volume(1h) = column E) with the latest values from the feed.
# Data rows live in Sheet1, row 2..51 (row 1 is the header).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores values as plain text (e.g. "27.119.85",
# "1.001") rather than numbers. Excel auto-converts numeric-looking
# strings assigned via .Value, so force these cells to Text format
# first (one cell at a time, so the format reliably "sticks" before
# the value is written) to keep them stored as text like the source.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Write the refreshed cell values (new coin ranking snapshot).
$ws.Range("D2").Value = '27.119.85'
$ws.Range("E2").Value = '  -2.22%  '
$ws.Range("D3").Value = '1.822.03'
$ws.Range("E3").Value = '  -1.42%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -1.29%  '
$ws.Range("D5").Value = '310.84'
$ws.Range("E5").Value = '  -2.93%  '
$ws.Range("E6").Value = '  -1.06%  '
$ws.Range("D7").Value = '0.4230'
$ws.Range("E7").Value = '  -1.95%  '
$ws.Range("D8").Value = '0.3673'
$ws.Range("E8").Value = '  -1.85%  '
$ws.Range("D9").Value = '0.07226'
$ws.Range("E9").Value = '  -1.88%  '
$ws.Range("D10").Value = '0.8458'
$ws.Range("E10").Value = '  -3.73%  '
$ws.Range("D11").Value = '20.90'
$ws.Range("E11").Value = '  -3.66%  '
$ws.Range("D12").Value = '1.816.07'
$ws.Range("E12").Value = '  -1.83%  '
$ws.Range("D13").Value = '6.644'
$ws.Range("E13").Value = '  -1.36%  '
$ws.Range("E14").Value = '  -0.71%  '
$ws.Range("D15").Value = '5.284'
$ws.Range("E15").Value = '  -3.12%  '
$ws.Range("D16").Value = '89.39'
$ws.Range("E16").Value = '  +1.26%  '
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  -1.35%  '
$ws.Range("D18").Value = '0.000008830'
$ws.Range("E18").Value = '  -1.91%  '
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  -0.96%  '
$ws.Range("D20").Value = '15.00'
$ws.Range("E20").Value = '  -3.09%  '
$ws.Range("D21").Value = '27.142.63'
$ws.Range("E21").Value = '  -2.18%  '
$ws.Range("D22").Value = '5.115'
$ws.Range("E22").Value = '  -2.13%  '
$ws.Range("D23").Value = '10.85'
$ws.Range("E23").Value = '  -2.37%  '
$ws.Range("D24").Value = '2.045.22'
$ws.Range("E24").Value = '  -1.94%  '
$ws.Range("E25").Value = '  -2.13%  '
$ws.Range("D26").Value = '151.89'
$ws.Range("E26").Value = '  -2.46%  '
$ws.Range("D27").Value = '2.252'
$ws.Range("E27").Value = '  +4.89%  '
$ws.Range("D28").Value = '18.34'
$ws.Range("E28").Value = '  -1.57%  '
$ws.Range("D29").Value = '5.209'
$ws.Range("E29").Value = '  -3.54%  '
$ws.Range("D30").Value = '116.18'
$ws.Range("E30").Value = '  -2.37%  '
$ws.Range("D31").Value = '0.08808'
$ws.Range("E31").Value = '  -1.73%  '
$ws.Range("D32").Value = '1.181'
$ws.Range("E32").Value = '  -4.27%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '0.7411'
$ws.Range("E33").Value = '  -4.94%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '2.959'
$ws.Range("E34").Value = '  +1.26%  '
$ws.Range("D35").Value = '4.423'
$ws.Range("E35").Value = '  -3.10%  '
$ws.Range("D36").Value = '1.000'
$ws.Range("E36").Value = '  -1.22%  '
$ws.Range("D37").Value = '1.098'
$ws.Range("E37").Value = '  -3.35%  '
$ws.Range("D38").Value = '0.01968'
$ws.Range("E38").Value = '  -0.25%  '
$ws.Range("D39").Value = '0.05243'
$ws.Range("E39").Value = '  -2.07%  '
$ws.Range("D40").Value = '7.300'
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("D41").Value = '2.871'
$ws.Range("E41").Value = '  -0.47%  '
$ws.Range("D42").Value = '0.1694'
$ws.Range("E42").Value = '  +0.29%  '
$ws.Range("D43").Value = '0.5032'
$ws.Range("E43").Value = '  -2.16%  '
$ws.Range("D44").Value = '8.588'
$ws.Range("E44").Value = '  -2.60%  '
$ws.Range("D45").Value = '10.59'
$ws.Range("E45").Value = '  -1.22%  '
$ws.Range("D46").Value = '0.4754'
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Value = '106.15'
$ws.Range("E47").Value = '  -3.27%  '
$ws.Range("E48").Value = '  -1.19%  '
$ws.Range("D49").Value = '0.06372'
$ws.Range("E49").Value = '  -1.99%  '
$ws.Range("D50").Value = '1.653'
$ws.Range("E50").Value = '  -2.54%  '
$ws.Range("D51").Value = '1.883'
$ws.Range("E51").Value = '  +1.29%  '
